# Generate Report for Handoff
# Update status strings from "In Translation" to "Ready for handoff" and
# refresh the related timestamp columns to reflect the handoff generation.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns (E, F) and generate-date column (G)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-17 04:52:49"

# zh-cn sheet: Status (C) and Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-17 04:52:44"

# de-de sheet: Status (C) and Latest Handoff Datetime (H)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-17 04:52:49"

# Widen the status columns to fit the longer "Ready for handoff" text.
$overview.Columns.Item(5).ColumnWidth = 17.2159881591797
$overview.Columns.Item(6).ColumnWidth = 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 17.2159881591797
